$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 78-79, pushing the existing data (old rows 78-144)
# down to rows 80-146.
$ws.Rows("78:79").Insert()

# Row 78: new record (Red Globe / Primera)
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = "2022-04-28"
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100109
$ws.Range("H78").Value = "Uva"
$ws.Range("I78").Value = 100109001
$ws.Range("J78").Value = "Uva"
$ws.Range("K78").Value = "Red Globe"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 150
$ws.Range("N78").Value = 8000
$ws.Range("O78").Value = 9000
$ws.Range("P78").Value = 8467
$ws.Range("Q78").Value = "$/bandeja 18 kilos"
$ws.Range("R78").Value = "Región de O'Higgins"
$ws.Range("S78").Value = 470
$ws.Range("T78").Value = 18

# Row 79: new record (Thompson seedless / Primera)
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = "2022-04-28"
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100109
$ws.Range("H79").Value = "Uva"
$ws.Range("I79").Value = 100109001
$ws.Range("J79").Value = "Uva"
$ws.Range("K79").Value = "Thompson seedless"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 130
$ws.Range("N79").Value = 10000
$ws.Range("O79").Value = 11000
$ws.Range("P79").Value = 10385
$ws.Range("Q79").Value = "$/bandeja 18 kilos"
$ws.Range("R79").Value = "Región de O'Higgins"
$ws.Range("S79").Value = 577
$ws.Range("T79").Value = 18
